$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "partner"
$ws.Cells.Item(2, 2).Value = "相手|あいて"
$ws.Cells.Item(3, 1).Value = "reception desk"
$ws.Cells.Item(3, 2).Value = "受付|うけつけ"
$ws.Cells.Item(4, 1).Value = "English conversation"
$ws.Cells.Item(4, 2).Value = "英会話|えいかいわ"
$ws.Cells.Item(5, 1).Value = "(someone's) daughter (polite)"
$ws.Cells.Item(5, 2).Value = "お嬢さん|おじょうさん"
$ws.Cells.Item(6, 1).Value = "household matters"
$ws.Cells.Item(6, 2).Value = "家事|かじ"
$ws.Cells.Item(7, 1).Value = "wind"
$ws.Cells.Item(7, 2).Value = "風|かぜ"
$ws.Cells.Item(8, 1).Value = "God"
$ws.Cells.Item(8, 2).Value = "神様|かみさま"
$ws.Cells.Item(9, 1).Value = "tree"
$ws.Cells.Item(9, 2).Value = "木|き"
$ws.Cells.Item(10, 1).Value = "junior member of a group"
$ws.Cells.Item(10, 2).Value = "後輩|こうはい"
$ws.Cells.Item(11, 1).Value = "monkey"
$ws.Cells.Item(11, 2).Value = "猿|さる"
$ws.Cells.Item(12, 1).Value = "freedom"
$ws.Cells.Item(12, 2).Value = "自由|じゆう"
$ws.Cells.Item(13, 1).Value = "cram school"
$ws.Cells.Item(13, 2).Value = "塾|じゅく"
$ws.Cells.Item(14, 1).Value = "document"
$ws.Cells.Item(14, 2).Value = "書類|しょるい"
$ws.Cells.Item(15, 1).Value = "senior member of a group"
$ws.Cells.Item(15, 2).Value = "先輩|せんぱい"
$ws.Cells.Item(16, 1).Value = "living alone"
$ws.Cells.Item(16, 2).Value = "一人暮らし|ひとりぐらし"
$ws.Cells.Item(17, 1).Value = "subordinate"
$ws.Cells.Item(17, 2).Value = "部下|ぶか"
$ws.Cells.Item(18, 1).Value = "review of a lesson"
$ws.Cells.Item(18, 2).Value = "復習|ふくしゅう"
$ws.Cells.Item(19, 1).Value = "project"
$ws.Cells.Item(19, 2).Value = "プロジェクト"
$ws.Cells.Item(20, 1).Value = "ball"
$ws.Cells.Item(20, 2).Value = "ボール"
$ws.Cells.Item(21, 1).Value = "waste (money)"
$ws.Cells.Item(21, 2).Value = "無駄遣い|むだづかい"
$ws.Cells.Item(22, 1).Value = "license"
$ws.Cells.Item(22, 2).Value = "免許|めんきょ"
$ws.Cells.Item(23, 1).Value = "Europe"
$ws.Cells.Item(23, 2).Value = "ヨーロッパ"
$ws.Cells.Item(24, 1).Value = "preparation of lessons"
$ws.Cells.Item(24, 2).Value = "予習|よしゅう"
$ws.Cells.Item(25, 1).Value = "noisy; annoying"
$ws.Cells.Item(25, 2).Value = "うるさい"
$ws.Cells.Item(26, 1).Value = "worried about"
$ws.Cells.Item(26, 2).Value = "心配|しんぱい（な）"
$ws.Cells.Item(27, 1).Value = "poor"
$ws.Cells.Item(27, 2).Value = "貧乏|びんぼう（な）"
$ws.Cells.Item(28, 1).Value = "fluent"
$ws.Cells.Item(28, 2).Value = "ぺらぺら（な）"
$ws.Cells.Item(29, 1).Value = "strange; unusual"
$ws.Cells.Item(29, 2).Value = "変|へん（な）"
$ws.Cells.Item(30, 1).Value = "easy; comfortable"
$ws.Cells.Item(30, 2).Value = "楽|らく（な）"
$ws.Cells.Item(31, 1).Value = "the wind blows"
$ws.Cells.Item(31, 2).Value = "風が吹く|かぜがふく"
$ws.Cells.Item(32, 1).Value = "to win"
$ws.Cells.Item(32, 2).Value = "勝つ|かつ"
$ws.Cells.Item(33, 1).Value = "to make a photocopy"
$ws.Cells.Item(33, 2).Value = "コピーを取る|コピーをとる"
$ws.Cells.Item(34, 1).Value = "to carry"
$ws.Cells.Item(34, 2).Value = "運ぶ|はこぶ"
$ws.Cells.Item(35, 1).Value = "to run"
$ws.Cells.Item(35, 2).Value = "走る|はしる"
$ws.Cells.Item(36, 1).Value = "to pick up (something)"
$ws.Cells.Item(36, 2).Value = "拾う|ひろう"
$ws.Cells.Item(37, 1).Value = "to leave (someone/something) alone; to neglect"
$ws.Cells.Item(37, 2).Value = "放っておく|ほうっておく"
$ws.Cells.Item(38, 1).Value = "to be in time"
$ws.Cells.Item(38, 2).Value = "間に合う|まにあう"
$ws.Cells.Item(39, 1).Value = "to make a plan"
$ws.Cells.Item(39, 2).Value = "計画を立てる|けいかくをたてる"
$ws.Cells.Item(40, 1).Value = "to raise; to bring up"
$ws.Cells.Item(40, 2).Value = "育てる|そだてる"
$ws.Cells.Item(41, 1).Value = "to help; to rescue"
$ws.Cells.Item(41, 2).Value = "助ける|たすける"
$ws.Cells.Item(42, 1).Value = "to lose (a match)"
$ws.Cells.Item(42, 2).Value = "負ける|まける"
$ws.Cells.Item(43, 1).Value = "to pray for help"
$ws.Cells.Item(43, 2).Value = "お願いする|おねがいする"
$ws.Cells.Item(44, 1).Value = "to agree"
$ws.Cells.Item(44, 2).Value = "賛成する|さんせいする"
$ws.Cells.Item(45, 1).Value = "to fail; to be unsuccessful"
$ws.Cells.Item(45, 2).Value = "失敗する|しっぱいする"
$ws.Cells.Item(46, 1).Value = "to stay up all night"
$ws.Cells.Item(46, 2).Value = "徹夜する|てつやする"
$ws.Cells.Item(47, 1).Value = "to oppose; to object to"
$ws.Cells.Item(47, 2).Value = "反対する|はんたいする"
$ws.Cells.Item(48, 1).Value = "to translate"
$ws.Cells.Item(48, 2).Value = "翻訳する|ほんやくする"
$ws.Cells.Item(57, 1).Value = "diary"
$ws.Cells.Item(57, 2).Value = "日記|にっき"
$ws.Cells.Item(58, 1).Value = "to fill in"
$ws.Cells.Item(58, 2).Value = "記入する|きにゅうする"
$ws.Cells.Item(59, 1).Value = "an article; news"
$ws.Cells.Item(59, 2).Value = "記事|きじ"
$ws.Cells.Item(60, 1).Value = "to memorize"
$ws.Cells.Item(60, 2).Value = "暗記する|あんきする"
$ws.Cells.Item(61, 1).Value = "bank"
$ws.Cells.Item(61, 2).Value = "銀行|ぎんこう"
$ws.Cells.Item(62, 1).Value = "silver medal"
$ws.Cells.Item(62, 2).Value = "銀メダル|ぎんメダル"
$ws.Cells.Item(63, 1).Value = "land covered with snow"
$ws.Cells.Item(63, 2).Value = "銀世界|ぎんせかい"
$ws.Cells.Item(64, 1).Value = "one time"
$ws.Cells.Item(64, 2).Value = "一回|いっかい"
$ws.Cells.Item(65, 1).Value = "out-of-service bus"
$ws.Cells.Item(65, 2).Value = "回送バス|かいそうバス"
$ws.Cells.Item(66, 1).Value = "last inning; last episode"
$ws.Cells.Item(66, 2).Value = "最終回|さいしゅうかい"
$ws.Cells.Item(67, 1).Value = "to turn"
$ws.Cells.Item(67, 2).Value = "回す|まわす"
$ws.Cells.Item(68, 1).Value = "evening"
$ws.Cells.Item(68, 2).Value = "夕方|ゆうがた"
$ws.Cells.Item(69, 1).Value = "dinner"
$ws.Cells.Item(69, 2).Value = "夕食|ゆうしょく"
$ws.Cells.Item(70, 1).Value = "Tanabata"
$ws.Cells.Item(70, 2).Value = "七夕|たなばた"
$ws.Cells.Item(71, 1).Value = "setting sun"
$ws.Cells.Item(71, 2).Value = "夕日|ゆうひ"
$ws.Cells.Item(72, 1).Value = "evening newspaper"
$ws.Cells.Item(72, 2).Value = "夕刊|ゆうかん"
$ws.Cells.Item(73, 1).Value = "Mr./Ms. Kuroki"
$ws.Cells.Item(73, 2).Value = "黒木さん|くろきさん"
$ws.Cells.Item(74, 1).Value = "black"
$ws.Cells.Item(74, 2).Value = "黒い|くろい"
$ws.Cells.Item(75, 1).Value = "black and white photograph"
$ws.Cells.Item(75, 2).Value = "白黒写真|しろくろしゃしん"
$ws.Cells.Item(76, 1).Value = "blackboard"
$ws.Cells.Item(76, 2).Value = "黒板|こくばん"
$ws.Cells.Item(77, 1).Value = "a thing to take care of"
$ws.Cells.Item(77, 2).Value = "用事|ようじ"
$ws.Cells.Item(78, 1).Value = "to prepare"
$ws.Cells.Item(78, 2).Value = "用意する|よういする"
$ws.Cells.Item(79, 1).Value = "for children"
$ws.Cells.Item(79, 2).Value = "子供用|こどもよう"
$ws.Cells.Item(80, 1).Value = "cost"
$ws.Cells.Item(80, 2).Value = "費用|ひよう"
$ws.Cells.Item(81, 1).Value = "absence; not at home"
$ws.Cells.Item(81, 2).Value = "留守|るす"
$ws.Cells.Item(82, 1).Value = "answering machine"
$ws.Cells.Item(82, 2).Value = "留守番電話|るすばんでんわ"
$ws.Cells.Item(83, 1).Value = "charm"
$ws.Cells.Item(83, 2).Value = "お守り|おまもり"
$ws.Cells.Item(84, 1).Value = "security guard"
$ws.Cells.Item(84, 2).Value = "守衛|しゅえい"
$ws.Cells.Item(85, 1).Value = "weekend"
$ws.Cells.Item(85, 2).Value = "週末|しゅうまつ"
$ws.Cells.Item(86, 1).Value = "end of the month"
$ws.Cells.Item(86, 2).Value = "月末|げつまつ"
$ws.Cells.Item(87, 1).Value = "year-end"
$ws.Cells.Item(87, 2).Value = "年末|ねんまつ"
$ws.Cells.Item(88, 1).Value = "final examination"
$ws.Cells.Item(88, 2).Value = "期末試験|きまつしけん"
$ws.Cells.Item(89, 1).Value = "the end"
$ws.Cells.Item(89, 2).Value = "末|すえ"
$ws.Cells.Item(90, 1).Value = "to wait"
$ws.Cells.Item(90, 2).Value = "待つ|まつ"
$ws.Cells.Item(91, 1).Value = "waiting room"
$ws.Cells.Item(91, 2).Value = "待合室|まちあいしつ"
$ws.Cells.Item(92, 1).Value = "to expect"
$ws.Cells.Item(92, 2).Value = "期待する|きたいする"
$ws.Cells.Item(93, 1).Value = "invitation"
$ws.Cells.Item(93, 2).Value = "招待|しょうたい"
$ws.Cells.Item(94, 1).Value = "over-time work"
$ws.Cells.Item(94, 2).Value = "残業|ざんぎょう"
$ws.Cells.Item(95, 1).Value = "to leave"
$ws.Cells.Item(95, 2).Value = "残す|のこす"
$ws.Cells.Item(96, 1).Value = "regrettable"
$ws.Cells.Item(96, 2).Value = "残念|ざんねん"
$ws.Cells.Item(97, 1).Value = "regret"
$ws.Cells.Item(97, 2).Value = "心残り|こころのこり"
$ws.Cells.Item(98, 1).Value = "account balance"
$ws.Cells.Item(98, 2).Value = "残高|ざんだか"
$ws.Cells.Item(99, 1).Value = "answering machine"
$ws.Cells.Item(99, 2).Value = "留守番電話|るすばんでんわ"
$ws.Cells.Item(100, 1).Value = "the first"
$ws.Cells.Item(100, 2).Value = "一番|いちばん"
$ws.Cells.Item(101, 1).Value = "number"
$ws.Cells.Item(101, 2).Value = "番号|ばんごう"
$ws.Cells.Item(102, 1).Value = "TV program"
$ws.Cells.Item(102, 2).Value = "番組|ばんぐみ"
$ws.Cells.Item(103, 1).Value = "station"
$ws.Cells.Item(103, 2).Value = "駅|えき"
$ws.Cells.Item(104, 1).Value = "Tokyo Station"
$ws.Cells.Item(104, 2).Value = "東京駅|とうきょうえき"
$ws.Cells.Item(105, 1).Value = "station attendant"
$ws.Cells.Item(105, 2).Value = "駅員|えきいん"
$ws.Cells.Item(106, 1).Value = "vicinity; in front of the station"
$ws.Cells.Item(106, 2).Value = "駅前|えきまえ"
$ws.Cells.Item(107, 1).Value = "to explain"
$ws.Cells.Item(107, 2).Value = "説明する|せつめいする"
$ws.Cells.Item(108, 1).Value = "novel"
$ws.Cells.Item(108, 2).Value = "小説|しょうせつ"
$ws.Cells.Item(109, 1).Value = "novelist"
$ws.Cells.Item(109, 2).Value = "小説家|しょうせつか"
$ws.Cells.Item(110, 1).Value = "to preach"
$ws.Cells.Item(110, 2).Value = "説教する|せっきょうする"
$ws.Cells.Item(111, 1).Value = "to guide"
$ws.Cells.Item(111, 2).Value = "案内する|あんないする"
$ws.Cells.Item(112, 1).Value = "information desk"
$ws.Cells.Item(112, 2).Value = "案内所|あんないじょ"
$ws.Cells.Item(113, 1).Value = "idea; proposal"
$ws.Cells.Item(113, 2).Value = "案|あん"
$ws.Cells.Item(114, 1).Value = "proposal"
$ws.Cells.Item(114, 2).Value = "提案|ていあん"
$ws.Cells.Item(115, 1).Value = "my wife"
$ws.Cells.Item(115, 2).Value = "家内|かない"
$ws.Cells.Item(116, 1).Value = "domestic"
$ws.Cells.Item(116, 2).Value = "国内|こくない"
$ws.Cells.Item(117, 1).Value = "internal medicine"
$ws.Cells.Item(117, 2).Value = "内科|ないか"
$ws.Cells.Item(118, 1).Value = "inside"
$ws.Cells.Item(118, 2).Value = "内側|うちがわ"
$ws.Cells.Item(119, 1).Value = "to forget"
$ws.Cells.Item(119, 2).Value = "忘れる|わすれる"
$ws.Cells.Item(120, 1).Value = "lost article"
$ws.Cells.Item(120, 2).Value = "忘れ物|わすれもの"
$ws.Cells.Item(121, 1).Value = "year-end party"
$ws.Cells.Item(121, 2).Value = "忘年会|ぼうねんかい"
